$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2643.6667
$ws.Range("J17").Value = 2717.5862
$ws.Range("L17").Value = 8152.758600000001
$ws.Range("N17").Value = -8488.758600000001

$ws.Range("H19").Value = 14286082
$ws.Range("I19").Value = 23809836
$ws.Range("J19").Value = 450
$ws.Range("K19").Value = 23809836
$ws.Range("L19").Value = 450
$ws.Range("M19").Value = -23809661
$ws.Range("N19").Value = -800

$ws.Range("H62").Value = 3473
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 3966.25
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 3966.25
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -5214.25

$ws.Range("H65").Value = 3473
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 3966.25
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 19831.25
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -26071.25

$ws.Range("H111").Value = 660.53845
$ws.Range("I111").Value = 590.25
$ws.Range("J111").Value = 773
$ws.Range("K111").Value = 1770.75
$ws.Range("L111").Value = 2319
$ws.Range("M111").Value = 1296.25
$ws.Range("N111").Value = -8453

$ws.Range("H125").Value = 1475
$ws.Range("I125").Value = 1442.8572
$ws.Range("J125").Value = 1520
$ws.Range("K125").Value = 12985.7148
$ws.Range("L125").Value = 13680
$ws.Range("M125").Value = -10525.7148
$ws.Range("N125").Value = -18600

$ws.Range("H132").Value = 5266024.5
$ws.Range("I132").Value = 6063191.5
$ws.Range("K132").Value = 18189574.5
$ws.Range("M132").Value = -18187044.5

$ws.Range("H137").Value = 4004913.2
$ws.Range("I137").Value = 5004916.5
$ws.Range("K137").Value = 15014749.5
$ws.Range("M137").Value = -15012199.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 37138
$ws.Range("I25").Value = 4258
$ws.Range("J25").Value = 70018
$ws.Range("K25").Value = 4258
$ws.Range("L25").Value = 70018
$ws.Range("M25").Value = -3856
$ws.Range("N25").Value = -70822

$ws.Range("H32").Value = 7387.8877
$ws.Range("I32").Value = 5639.3887
$ws.Range("J32").Value = 23124.375
$ws.Range("K32").Value = 5639.3887
$ws.Range("L32").Value = 23124.375
$ws.Range("M32").Value = -5352.3887
$ws.Range("N32").Value = -23698.375

$ws.Range("H45").Value = 1481.0625
$ws.Range("I45").Value = 1090.8695
$ws.Range("J45").Value = 2478.2222
$ws.Range("K45").Value = 1090.8695
$ws.Range("L45").Value = 2478.2222
$ws.Range("M45").Value = -713.8695
$ws.Range("N45").Value = -3232.2222

$ws.Range("H61").Value = 2807.7778
$ws.Range("I61").Value = 1809.6923
$ws.Range("K61").Value = 1809.6923
$ws.Range("M61").Value = -1597.6923

$ws.Range("H74").Value = 1439.1111
$ws.Range("I74").Value = 1575.1818
$ws.Range("J74").Value = 1225.2858
$ws.Range("K74").Value = 1575.1818
$ws.Range("L74").Value = 1225.2858
$ws.Range("M74").Value = -701.1818000000001
$ws.Range("N74").Value = -2973.2858

$ws.Range("H77").Value = 1439.1111
$ws.Range("I77").Value = 1575.1818
$ws.Range("J77").Value = 1225.2858
$ws.Range("K77").Value = 7875.909000000001
$ws.Range("L77").Value = 6126.429
$ws.Range("M77").Value = -3507.909000000001
$ws.Range("N77").Value = -14862.429

$ws.Range("H102").Value = 2797
$ws.Range("I102").Value = 2261.4
$ws.Range("J102").Value = 5475
$ws.Range("K102").Value = 2261.4
$ws.Range("L102").Value = 5475
$ws.Range("M102").Value = -639.4000000000001
$ws.Range("N102").Value = -8719

$ws.Range("H136").Value = 2807.7778
$ws.Range("I136").Value = 1809.6923
$ws.Range("K136").Value = 5429.0769
$ws.Range("M136").Value = -2879.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 667.2222
$ws.Range("I64").Value = 549.5
$ws.Range("J64").Value = 761.4
$ws.Range("K64").Value = 549.5
$ws.Range("L64").Value = 761.4
$ws.Range("M64").Value = -324.5
$ws.Range("N64").Value = -1211.4

$ws.Range("H67").Value = 667.2222
$ws.Range("I67").Value = 549.5
$ws.Range("J67").Value = 761.4
$ws.Range("K67").Value = 549.5
$ws.Range("L67").Value = 761.4
$ws.Range("M67").Value = 230.5
$ws.Range("N67").Value = -2321.4

$ws.Range("H105").Value = 1847.1666
$ws.Range("I105").Value = 1669.9333
$ws.Range("J105").Value = 2733.3333
$ws.Range("K105").Value = 1669.9333
$ws.Range("L105").Value = 2733.3333
$ws.Range("M105").Value = 77.06670000000008
$ws.Range("N105").Value = -6227.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1589242.6
$ws.Range("I31").Value = 1888248.2
$ws.Range("K31").Value = 1888248.2
$ws.Range("M31").Value = -1887953.2

$ws.Range("H34").Value = 1589242.6
$ws.Range("I34").Value = 1888248.2
$ws.Range("K34").Value = 1888248.2
$ws.Range("M34").Value = -1888046.2

$ws.Range("H58").Value = 11908319
$ws.Range("I58").Value = 2038.8462
$ws.Range("J58").Value = 31256024
$ws.Range("K58").Value = 2038.8462
$ws.Range("L58").Value = 31256024
$ws.Range("M58").Value = -1835.8462
$ws.Range("N58").Value = -31256430

$ws.Range("H94").Value = 1673.8788
$ws.Range("I94").Value = 1651.2
$ws.Range("J94").Value = 1677.9286
$ws.Range("K94").Value = 1651.2
$ws.Range("L94").Value = 1677.9286
$ws.Range("M94").Value = -1200.2
$ws.Range("N94").Value = -2579.9286

$ws.Range("H99").Value = 2000.28
$ws.Range("I99").Value = 1610.3334
$ws.Range("K99").Value = 1610.3334
$ws.Range("M99").Value = -112.3334

$ws.Range("H105").Value = 4661
$ws.Range("I105").Value = 5377.5
$ws.Range("K105").Value = 5377.5
$ws.Range("M105").Value = -3630.5

$ws.Range("H122").Value = 1686.1
$ws.Range("I122").Value = 1417.4445
$ws.Range("J122").Value = 2244.077
$ws.Range("K122").Value = 4252.333500000001
$ws.Range("L122").Value = 6732.231000000001
$ws.Range("M122").Value = -1802.333500000001
$ws.Range("N122").Value = -11632.231

$ws.Range("H126").Value = 2000.28
$ws.Range("I126").Value = 1610.3334
$ws.Range("K126").Value = 4831.0002
$ws.Range("M126").Value = -2361.0002

$ws.Range("H132").Value = 2314.611
$ws.Range("I132").Value = 2025.0476
$ws.Range("J132").Value = 2720
$ws.Range("K132").Value = 6075.142800000001
$ws.Range("L132").Value = 8160
$ws.Range("M132").Value = -3545.142800000001
$ws.Range("N132").Value = -13220

$ws.Range("H134").Value = 5054.3335
$ws.Range("I134").Value = 3372.25
$ws.Range("J134").Value = 6400
$ws.Range("K134").Value = 10116.75
$ws.Range("L134").Value = 19200
$ws.Range("M134").Value = -7581.75
$ws.Range("N134").Value = -24270

$ws.Range("H136").Value = 11908319
$ws.Range("I136").Value = 2038.8462
$ws.Range("J136").Value = 31256024
$ws.Range("K136").Value = 6116.5386
$ws.Range("L136").Value = 93768072
$ws.Range("M136").Value = -3566.5386
$ws.Range("N136").Value = -93773172

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 198.13333
$ws.Range("I98").Value = 198.66667
$ws.Range("J98").Value = 196
$ws.Range("K98").Value = 596.00001
$ws.Range("L98").Value = 588
$ws.Range("M98").Value = 901.99999
$ws.Range("N98").Value = -3584

$ws.Range("H113").Value = 2381719.8
$ws.Range("I113").Value = 7692847.5
$ws.Range("J113").Value = 869.2069
$ws.Range("K113").Value = 23078542.5
$ws.Range("L113").Value = 2607.6207
$ws.Range("M113").Value = -23076372.5
$ws.Range("N113").Value = -6947.620699999999

$ws.Range("H131").Value = 2836.4707
$ws.Range("J131").Value = 3300
$ws.Range("L131").Value = 9900
$ws.Range("N131").Value = -19980

$ws.Range("H133").Value = 5370.909
$ws.Range("I133").Value = 5912.857
$ws.Range("J133").Value = 4422.5
$ws.Range("K133").Value = 17738.571
$ws.Range("L133").Value = 13267.5
$ws.Range("M133").Value = -12678.571
$ws.Range("N133").Value = -23387.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 50807.2
$ws.Range("I22").Value = 9000
$ws.Range("J22").Value = 61259
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 61259
$ws.Range("M22").Value = -8471
$ws.Range("N22").Value = -62317

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

$ws.Range("H103").Value = 39920
$ws.Range("J103").Value = 39920
$ws.Range("L103").Value = 39920
$ws.Range("N103").Value = -42264

$ws.Range("H132").Value = 2463.7942
$ws.Range("J132").Value = 3363.4119
$ws.Range("L132").Value = 10090.2357
$ws.Range("N132").Value = -15150.2357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 8599.666999999999
$ws.Range("J19").Value = 8599.666999999999
$ws.Range("L19").Value = 8599.666999999999
$ws.Range("N19").Value = -8947.666999999999

$ws.Range("H136").Value = 1328.2413
$ws.Range("I136").Value = 630.9048
$ws.Range("J136").Value = 3158.75
$ws.Range("K136").Value = 1892.7144
$ws.Range("L136").Value = 9476.25
$ws.Range("M136").Value = 657.2855999999999
$ws.Range("N136").Value = -14576.25
